# ---------------------------------------------------------------------------
# Applies the commit:
#   "added linux postinstall, fixed floating points, fixed default data,
#    added classic save features"
#
#  * Reset/refresh the legacy 56-slot colour palette (classic-save artifact
#    the re-exported workbook now carries in styles.xml).
#  * "fixed floating points" / "fixed default data": the Grade column on the
#    Students & Grades sheets used zero-padded strings ("07","02",...); they
#    are normalised to their plain form ("7","2",...). Student E's Task 1
#    score/point total is corrected (0 -> 1 point, 3 (9.7%) -> 4 (12.9%)).
#  * Two new summary sheets, "Points" and "Score", are appended that
#    aggregate the Students sheet by point total / by letter grade.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- "added classic save features": refresh the legacy colour palette -----
try { $wb.ResetColors() } catch { }

# --- "fixed floating points" / "fixed default data" ------------------------
$students = $wb.Worksheets.Item("Students")
$students.Range("B4").Value  = "7"
$students.Range("B5").Value  = "2"
$students.Range("B6").Value  = "0"
$students.Range("D6").Value  = "4 (12.9%)"
$students.Range("E6").Value  = 1
$students.Range("B8").Value  = "4"
$students.Range("B9").Value  = "8"
$students.Range("B10").Value = "9"

$grades = $wb.Worksheets.Item("Grades")
$grades.Range("A8").Value  = "9"
$grades.Range("A9").Value  = "8"
$grades.Range("A10").Value = "7"
$grades.Range("A11").Value = "6"
$grades.Range("A12").Value = "5"
$grades.Range("A13").Value = "4"
$grades.Range("A14").Value = "3"
$grades.Range("A15").Value = "2"
$grades.Range("A16").Value = "1"
$grades.Range("A17").Value = "0"

# --- add the two new summary sheets, at the end of the tab strip -----------
$additional = $wb.Worksheets.Item("Additional")
$points = $wb.Worksheets.Add([System.Type]::Missing, $additional)
$points.Name = "Points"
$score  = $wb.Worksheets.Add([System.Type]::Missing, $points)
$score.Name  = "Score"

# Reuse the bold/centered/bordered header style already used by the other
# sheets (cellXfs index 1) instead of minting a new one.
$students.Range("A1:D1").Copy()
$points.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats
$students.Range("A1:C1").Copy()
$score.Range("A1:C1").PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Points sheet: per distinct point-total summary of the Students sheet
# ---------------------------------------------------------------------------
$points.Range("A1").Value = "Points"
$points.Range("B1").Value = "Grade"
$points.Range("C1").Value = "Amount"
$points.Range("D1").Value = "Students"

$pointsRows = @(
    @(3,  "0", 1, "Student E"),
    @(9,  "2", 1, "Student D"),
    @(13, "4", 1, "Student G"),
    @(17, "7", 1, "Student C"),
    @(19, "8", 1, "Student H"),
    @(20, "9", 1, "Student I"),
    @(22, "10", 2, "Student A, Student J"),
    @(27, "13", 2, "Student B, Student F")
)

$points.Range("B2:B9").NumberFormat = "@"
$r = 2
foreach ($row in $pointsRows) {
    $points.Cells.Item($r, 1).Value = $row[0]
    $points.Cells.Item($r, 2).Value = $row[1]
    $points.Cells.Item($r, 3).Value = $row[2]
    $points.Cells.Item($r, 4).Value = $row[3]
    $r++
}
$points.Range("A2:A9").Copy()
$points.Range("A2:A9").PasteSpecial(-4122)   # xlPasteFormats (border column)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Score sheet: per letter-grade summary of the Students sheet
# ---------------------------------------------------------------------------
$score.Range("A1").Value = "Grade"
$score.Range("B1").Value = "Amount"
$score.Range("C1").Value = "Students"

$scoreRows = @(
    @("0",  1, "Student E"),
    @("10", 2, "Student A, Student J"),
    @("13", 2, "Student B, Student F"),
    @("2",  1, "Student D"),
    @("4",  1, "Student G"),
    @("7",  1, "Student C"),
    @("8",  1, "Student H"),
    @("9",  1, "Student I")
)

$score.Range("A2:A9").NumberFormat = "@"
$r = 2
foreach ($row in $scoreRows) {
    $score.Cells.Item($r, 1).Value = $row[0]
    $score.Cells.Item($r, 2).Value = $row[1]
    $score.Cells.Item($r, 3).Value = $row[2]
    $r++
}

Write-Output "applied"
